# Insert a new data row before the existing row 135. This pushes the
# existing rows 135-212 down to 136-213 (matching the dimension change
# from A1:R212 to A1:R213) and keeps the row-above formatting (style
# index 2, the date/time number format already used by column D).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("135").Insert()

# Fill in the values for the newly inserted row 135.
$ws.Range("A135").Value = 10
$ws.Range("B135").Value = "Vega Modelo de Temuco"
$ws.Range("C135").Value = "La Araucanía"
$ws.Range("D135").Value = 44603
$ws.Range("E135").Value = 9
$ws.Range("F135").Value = 100112039
$ws.Range("G135").Value = "Ciboulette"
$ws.Range("H135").Value = "Sin especificar"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 65
$ws.Range("K135").Value = 5000
$ws.Range("L135").Value = 5000
$ws.Range("M135").Value = 5000
$ws.Range("N135").Value = "`$/docena de atados"
$ws.Range("O135").Value = "Provincia de Cautín"
$ws.Range("P135").Value = 1667
$ws.Range("Q135").Value = 3
$ws.Range("R135").Value = "Hortaliza"
